# B6-PowerPoint.pptx edit:
#   1) Re-style the three summary tables (slides 14, 15, 16) from the
#      deck's custom "Table_0" style to the built-in PowerPoint table
#      style {0157FE65-9905-422E-A318-A886696800E6}.
#   2) Swap the presentation's colour theme from "Integral" (Red Violet)
#      to the standard "Office Theme" colours.

$p = $ppt.ActivePresentation

# --- 1) Table styles -------------------------------------------------
$newTableStyleId = "{0157FE65-9905-422E-A318-A886696800E6}"
$tableSlideIndexes = @(14, 15, 16)

foreach ($slideIdx in $tableSlideIndexes) {
    $slide = $p.Slides.Item($slideIdx)
    for ($shapeIdx = 1; $shapeIdx -le $slide.Shapes.Count; $shapeIdx++) {
        $shape = $slide.Shapes.Item($shapeIdx)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# --- 2) Theme colours --------------------------------------------------
# Office Theme colour values (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
$officeThemeRGB = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = $officeThemeRGB[$i - 1]
}
